$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - new entry: Toblerone
$ws.Range("E21").Value = "Toblerone"
$ws.Range("C21").Value = "00:05:34"
$ws.Range("D21").Value = "00:09:16"
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 6.5

# Row 20 - new entry: kinder Riegel
$ws.Range("C20").Value = "00:29:19"
$ws.Range("E20").Value = "kinder Riegel"
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 6

# Correct candy-bar naming (rows 13 and 15)
$ws.Range("E13").Value = "Milky Way"
$ws.Range("E15").Value = "hanuta"

$ws.Range("E29").Select()
